# Update the "Kode_PKS" keyword value in cell A2 (sheet "inquiry")
# from "01733709" to "01733710", keeping it a text value (leading
# apostrophe preserves the existing quote-prefixed text formatting).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "'01733710"
